$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
